$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "27.01.2025"
$ws.Range("A3").Value = "27.01.2025"
